$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Hola muchachos"
$ws.Range("D2").Value = "Ronfis"

$ws.Range("E6").Select()
